$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.74808933333333
$ws.Range("H2").Value = 47.244268
$ws.Range("I2").Value = 0.2237852983702856
$ws.Range("J2").Value = 0.2380674495116221
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 218.8881844778142
$ws.Range("R2").Value = 1969.993660300328
$ws.Range("S2").Value = 0.01086143635920567
$ws.Range("T2").Value = 0.01190164381638466

$ws.Range("G3").Value = 15.74808933333333
$ws.Range("H3").Value = 47.244268
$ws.Range("I3").Value = 0.2237852983702856
$ws.Range("J3").Value = 0.2380674495116221
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 1108.072599750082
$ws.Range("R3").Value = 9972.653397750739
$ws.Range("S3").Value = 0.05498359837136366
$ws.Range("T3").Value = 0.06024941655202735

$ws.Range("G4").Value = 15.74808933333333
$ws.Range("H4").Value = 47.244268
$ws.Range("I4").Value = 0.2237852983702856
$ws.Range("J4").Value = 0.2380674495116221
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 1291.464034438073
$ws.Range("R4").Value = 11623.17630994266
$ws.Range("S4").Value = 0.06408365281897561
$ws.Range("T4").Value = 0.0702209896629253

$ws.Range("G5").Value = 15.74808933333333
$ws.Range("H5").Value = 47.244268
$ws.Range("I5").Value = 0.2237852983702856
$ws.Range("J5").Value = 0.2380674495116221
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 394.4919448950873
$ws.Range("R5").Value = 2366.951669370524
$ws.Range("S5").Value = 0.01957505912856411
$ws.Range("T5").Value = 0.01429985094223623

$ws.Range("G6").Value = 15.74808933333333
$ws.Range("H6").Value = 47.244268
$ws.Range("I6").Value = 0.2237852983702856
$ws.Range("J6").Value = 0.2380674495116221
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 1496.980091064542
$ws.Range("R6").Value = 13472.82081958088
$ws.Range("S6").Value = 0.07428155169217655
$ws.Range("T6").Value = 0.08139554853804855

$ws.Range("I7").Value = 0.366480229475165
$ws.Range("J7").Value = 0.38986928168634
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 358.4605095197043
$ws.Range("R7").Value = 3226.144585677338
$ws.Range("S7").Value = 0.0177871456183206
$ws.Range("T7").Value = 0.01949063315921329

$ws.Range("I8").Value = 0.366480229475165
$ws.Range("J8").Value = 0.38986928168634
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.09004345636309793
$ws.Range("T8").Value = 0.09866698198912435

$ws.Range("I9").Value = 0.366480229475165
$ws.Range("J9").Value = 0.38986928168634
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 2114.955893646998
$ws.Range("R9").Value = 19034.60304282299
$ws.Range("S9").Value = 0.1049460887812432
$ws.Range("T9").Value = 0.1149968500748444

$ws.Range("I10").Value = 0.366480229475165
$ws.Range("J10").Value = 0.38986928168634
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 646.0366232460798
$ws.Range("R10").Value = 3876.219739476479
$ws.Range("S10").Value = 0.0320569412453354
$ws.Range("T10").Value = 0.02341803815056709

$ws.Range("I11").Value = 0.366480229475165
$ws.Range("J11").Value = 0.38986928168634
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 2451.517643421442
$ws.Range("R11").Value = 22063.65879079298
$ws.Range("S11").Value = 0.1216465974671679
$ws.Range("T11").Value = 0.1332967783125909

$ws.Range("G12").Value = 2.416095
$ws.Range("H12").Value = 7.248285
$ws.Range("I12").Value = 0.03433346922420018
$ws.Range("J12").Value = 0.03652465783327086
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 33.58214681679
$ws.Range("R12").Value = 302.23932135111
$ws.Range("S12").Value = 0.001666377522049555
$ws.Range("T12").Value = 0.001825967678230165

$ws.Range("G13").Value = 2.416095
$ws.Range("H13").Value = 7.248285
$ws.Range("I13").Value = 0.03433346922420018
$ws.Range("J13").Value = 0.03652465783327086
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 170.002126050075
$ws.Range("R13").Value = 1530.019134450675
$ws.Range("S13").Value = 0.008435664434914721
$ws.Range("T13").Value = 0.009243553995858536

$ws.Range("G14").Value = 2.416095
$ws.Range("H14").Value = 7.248285
$ws.Range("I14").Value = 0.03433346922420018
$ws.Range("J14").Value = 0.03652465783327086
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 198.138309368175
$ws.Range("R14").Value = 1783.244784313575
$ws.Range("S14").Value = 0.009831808156557502
$ws.Range("T14").Value = 0.01077340739111328

$ws.Range("G15").Value = 2.416095
$ws.Range("H15").Value = 7.248285
$ws.Range("I15").Value = 0.03433346922420018
$ws.Range("J15").Value = 0.03652465783327086
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 60.5235337079175
$ws.Range("R15").Value = 363.141202247505
$ws.Range("S15").Value = 0.003003234327933375
$ws.Range("T15").Value = 0.002193904138526323

$ws.Range("G16").Value = 2.416095
$ws.Range("H16").Value = 7.248285
$ws.Range("I16").Value = 0.03433346922420018
$ws.Range("J16").Value = 0.03652465783327086
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 229.6688846859
$ws.Range("R16").Value = 2067.0199621731
$ws.Range("S16").Value = 0.01139638478274503
$ws.Range("T16").Value = 0.01248782462954256

$ws.Range("G17").Value = 12.6651745
$ws.Range("H17").Value = 25.330349
$ws.Range("I17").Value = 0.1799761097617747
$ws.Range("J17").Value = 0.1276415496938013
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 176.0376762996757
$ws.Range("R17").Value = 1056.226057798054
$ws.Range("S17").Value = 0.008735154081124793
$ws.Range("T17").Value = 0.006381150651814846

$ws.Range("G18").Value = 12.6651745
$ws.Range("H18").Value = 25.330349
$ws.Range("I18").Value = 0.1799761097617747
$ws.Range("J18").Value = 0.1276415496938013
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 891.1514620886991
$ws.Range("R18").Value = 5346.908772532194
$ws.Range("S18").Value = 0.04421976871424296
$ws.Range("T18").Value = 0.03230315153383748

$ws.Range("G19").Value = 12.6651745
$ws.Range("H19").Value = 25.330349
$ws.Range("I19").Value = 0.1799761097617747
$ws.Range("J19").Value = 0.1276415496938013
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 1038.641387562542
$ws.Range("R19").Value = 6231.848325375255
$ws.Range("S19").Value = 0.05153835670920394
$ws.Range("T19").Value = 0.03764948110292005

$ws.Range("G20").Value = 12.6651745
$ws.Range("H20").Value = 25.330349
$ws.Range("I20").Value = 0.1799761097617747
$ws.Range("J20").Value = 0.1276415496938013
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 317.2644766730642
$ws.Range("R20").Value = 1269.057906692257
$ws.Range("S20").Value = 0.01574295995300947
$ws.Range("T20").Value = 0.007666966392935169

$ws.Range("G21").Value = 12.6651745
$ws.Range("H21").Value = 25.330349
$ws.Range("I21").Value = 0.1799761097617747
$ws.Range("J21").Value = 0.1276415496938013
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 1203.924722234557
$ws.Range("R21").Value = 7223.548333407341
$ws.Range("S21").Value = 0.05973987030419349
$ws.Range("T21").Value = 0.04364080001229376

$ws.Range("G22").Value = 13.75232733333333
$ws.Range("H22").Value = 41.25698199999999
$ws.Range("I22").Value = 0.1954248931685745
$ws.Range("J22").Value = 0.2078970612749657
$ws.Range("M22").Value = 13.89934866666667
$ws.Range("N22").Value = 41.69804600000001
$ws.Range("O22").Value = 0.04853507553134179
$ws.Range("P22").Value = 0.04999273878390351
$ws.Range("Q22").Value = 191.1483925841302
$ws.Range("R22").Value = 1720.335533257172
$ws.Range("S22").Value = 0.009484961950641161
$ws.Range("T22").Value = 0.01039334347826054

$ws.Range("G23").Value = 13.75232733333333
$ws.Range("H23").Value = 41.25698199999999
$ws.Range("I23").Value = 0.1954248931685745
$ws.Range("J23").Value = 0.2078970612749657
$ws.Range("O23").Value = 0.245697991654417
$ws.Range("P23").Value = 0.253077086664408
$ws.Range("Q23").Value = 967.6460920631122
$ws.Range("R23").Value = 8708.814828568009
$ws.Range("S23").Value = 0.04801550377079775
$ws.Range("T23").Value = 0.05261398259356022

$ws.Range("G24").Value = 13.75232733333333
$ws.Range("H24").Value = 41.25698199999999
$ws.Range("I24").Value = 0.1954248931685745
$ws.Range("J24").Value = 0.2078970612749657
$ws.Range("M24").Value = 82.007665
$ws.Range("N24").Value = 246.022995
$ws.Range("O24").Value = 0.2863622109480123
$ws.Range("P24").Value = 0.2949625822722868
$ws.Range("Q24").Value = 1127.796252922343
$ws.Range("R24").Value = 10150.16627630109
$ws.Range("S24").Value = 0.0559623044820321
$ws.Range("T24").Value = 0.0613218540404837

$ws.Range("G25").Value = 13.75232733333333
$ws.Range("H25").Value = 41.25698199999999
$ws.Range("I25").Value = 0.1954248931685745
$ws.Range("J25").Value = 0.2078970612749657
$ws.Range("M25").Value = 25.0501465
$ws.Range("N25").Value = 50.100293
$ws.Range("O25").Value = 0.0874724982879541
$ws.Range("P25").Value = 0.06006638442832619
$ws.Range("Q25").Value = 344.4978144159543
$ws.Range("R25").Value = 2066.986886495726
$ws.Range("S25").Value = 0.01709430363311175
$ws.Range("T25").Value = 0.01248762480406137

$ws.Range("G26").Value = 13.75232733333333
$ws.Range("H26").Value = 41.25698199999999
$ws.Range("I26").Value = 0.1954248931685745
$ws.Range("J26").Value = 0.2078970612749657
$ws.Range("M26").Value = 95.05788666666668
$ws.Range("N26").Value = 285.17366
$ws.Range("O26").Value = 0.3319322235782747
$ws.Range("P26").Value = 0.3419012078510756
$ws.Range("Q26").Value = 1307.267173054902
$ws.Range("R26").Value = 11765.40455749412
$ws.Range("S26").Value = 0.06486781933199171
$ws.Range("T26").Value = 0.07108025635859984
